# Update on 2019-09-19, 更新学费，UI培训费用
#
# Target sheet: 第三学年 (3rd worksheet, "Worksheets.Item(3)") -- the active
# bill-ledger sheet. Rows 10-14 get new/changed entries:
#   - a new "第二学年剩余生活费" income row is inserted at the top (row 10),
#     pushing the previously-first row ("支出"/"生活费" on 2019-08-22) down to
#     row 11 (and its highlighted banding moves with it),
#   - two new 学费 (tuition) rows (收入 + 支出) are added for 2019-09-10,
#   - a new 培训费 (training fee) row is added for 2019-09-19.
# The summary formulas in row 4 (C4/E4/F4/G4) recalc automatically.

$wb = $excel.ActiveWorkbook

$wsYear3 = $wb.Worksheets.Item(3)   # 第三学年
$wsSheet1 = $wb.Worksheets.Item(4)  # Sheet1

# ------------------------------------------------------------------
# 1) Sheet1's stored (inactive) selection moves to G16 -- visit it first
#    so it keeps its own remembered cursor position without becoming the
#    active tab in the saved workbook.
# ------------------------------------------------------------------
$wsSheet1.Activate()
$wsSheet1.Range("G16").Select()

# ------------------------------------------------------------------
# 2) Switch to 第三学年 (stays the active/selected tab on save).
# ------------------------------------------------------------------
$wsYear3.Activate()

# ------------------------------------------------------------------
# 3) Row 10's shaded banding needs to move down to row 11 (the new row
#    10 reuses the format that row 11 used to have, and vice versa).
#    Stash row 10's current formatting off to one side first, then swap.
# ------------------------------------------------------------------
$wsYear3.Range("B10:G10").Copy()
$wsYear3.Range("Z200:AE200").PasteSpecial(-4122)

$wsYear3.Range("B11:G11").Copy()
$wsYear3.Range("B10:G10").PasteSpecial(-4122)

$wsYear3.Range("Z200:AE200").Copy()
$wsYear3.Range("B11:G11").PasteSpecial(-4122)

$wsYear3.Range("Z200:AE200").Clear()

# ------------------------------------------------------------------
# 4) Write the new ledger content for rows 10-14.
# ------------------------------------------------------------------

# Row 10: 第二学年剩余生活费 carried in as income
$wsYear3.Range("B10").Value = 1
$wsYear3.Range("C10").Value = "收入"
$wsYear3.Range("D10").Value = 9240
$wsYear3.Range("E10").Value = 43670
$wsYear3.Range("F10").Value = "生活费"
$wsYear3.Range("G10").Value = "第二学年剩余生活费"

# Row 11: previous first entry, now pushed down
$wsYear3.Range("B11").Value = 2
$wsYear3.Range("C11").Value = "支出"
$wsYear3.Range("D11").Value = 1200
$wsYear3.Range("E11").Value = 43699
$wsYear3.Range("F11").Value = "生活费"
$wsYear3.Range("G11").Value = "生活费"

# Row 12: tuition income
$wsYear3.Range("B12").Value = 3
$wsYear3.Range("C12").Value = "收入"
$wsYear3.Range("D12").Value = 16120
$wsYear3.Range("E12").Value = 43718
$wsYear3.Range("F12").Value = "学费"
$wsYear3.Range("G12").Value = "学费"

# Row 13: tuition expense (same amount/date, paired with row 12)
$wsYear3.Range("B13").Value = 4
$wsYear3.Range("C13").Value = "支出"
$wsYear3.Range("D13").Value = 16120
$wsYear3.Range("E13").Value = 43718
$wsYear3.Range("F13").Value = "学费"
$wsYear3.Range("G13").Value = "学费+住宿费等"

# Row 14: UI training-class fee
$wsYear3.Range("B14").Value = 5
$wsYear3.Range("C14").Value = "支出"
$wsYear3.Range("D14").Value = 1000
$wsYear3.Range("E14").Value = 43727
$wsYear3.Range("F14").Value = "培训费"
$wsYear3.Range("G14").Value = "UI培训班第一次培训费"

# ------------------------------------------------------------------
# 5) Recalculate the summary formulas (row 4) and leave the sheet's
#    selection where the author left it (G14).
# ------------------------------------------------------------------
$wb.Application.Calculate()
$wsYear3.Range("G14").Select()
